$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Facility demographics data update - add tri proximity tables
# Rewrite rows 2-11 with corrected/expanded data (new Chemours Chambers Works row inserted at row 4)

# Row 2: Daikin America Inc.
$ws.Range("A2").Value = "Daikin America Inc."
$ws.Range("B2").Value = "Decatur"
$ws.Range("C2").Value = 96067
$ws.Range("D2").Value = 66
$ws.Range("E2").Value = 552.179862807553
$ws.Range("F2").Value = 105557
$ws.Range("G2").Value = 191.164160647396
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0.41
$ws.Range("J2").Value = 77858
$ws.Range("K2").Value = 18646
$ws.Range("L2").Value = 1234
$ws.Range("M2").Value = 651
$ws.Range("N2").Value = 9961
$ws.Range("O2").Value = 53.8297692307692
$ws.Range("P2").Value = 6.23020145278463
$ws.Range("Q2").Value = 11.8536876343167
$ws.Range("R2").Value = 36.2121212121212
$ws.Range("S2").Value = 0.454545454545455

# Row 3: Chemours El Dorado
$ws.Range("A3").Value = "Chemours El Dorado"
$ws.Range("B3").Value = "El Dorado"
$ws.Range("C3").Value = 66990
$ws.Range("D3").Value = 37
$ws.Range("E3").Value = 1016.4657349007
$ws.Range("F3").Value = 36176
$ws.Range("G3").Value = 35.5899847460516
$ws.Range("H3").Value = 1
$ws.Range("I3").Value = 0.65
$ws.Range("J3").Value = 22446
$ws.Range("K3").Value = 12500
$ws.Range("L3").Value = 259
$ws.Range("M3").Value = 205
$ws.Range("N3").Value = 1621
$ws.Range("O3").Value = 47.55665625
$ws.Range("P3").Value = 9.05118823067004
$ws.Range("Q3").Value = 13.0819257713673
$ws.Range("R3").Value = 49.1891891891892
$ws.Range("S3").Value = 0.537837837837838

# Row 4: Chemours Chambers Works
$ws.Range("A4").Value = "Chemours Chambers Works"
$ws.Range("B4").Value = "Deepwater"
$ws.Range("C4").Value = 2619
$ws.Range("D4").Value = 306
$ws.Range("E4").Value = 465.982576285131
$ws.Range("F4").Value = 394757
$ws.Range("G4").Value = 847.149700632693
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0.12
$ws.Range("J4").Value = 251280
$ws.Range("K4").Value = 107713
$ws.Range("L4").Value = 1195
$ws.Range("M4").Value = 15665
$ws.Range("N4").Value = 42304
$ws.Range("O4").Value = 70.4335463917526
$ws.Range("P4").Value = 6.44340495189268
$ws.Range("Q4").Value = 7.50338941192498
$ws.Range("R4").Value = 32.2368421052632
$ws.Range("S4").Value = 0.347039473684211

# Row 5: ARKEMA, INC.
$ws.Range("A5").Value = "ARKEMA, INC."
$ws.Range("B5").Value = "Calvert City"
$ws.Range("C5").Value = 843010
$ws.Range("D5").Value = 29
$ws.Range("E5").Value = 698.358246272251
$ws.Range("F5").Value = 34960
$ws.Range("G5").Value = 50.0602666133208
$ws.Range("H5").Value = 1
$ws.Range("I5").Value = 0.83
$ws.Range("J5").Value = 33695
$ws.Range("K5").Value = 625
$ws.Range("L5").Value = 68
$ws.Range("M5").Value = 69
$ws.Range("N5").Value = 709
$ws.Range("O5").Value = 53.2628214285714
$ws.Range("P5").Value = 6.25440564467844
$ws.Range("Q5").Value = 5.87057212005925
$ws.Range("R5").Value = 31.4285714285714
$ws.Range("S5").Value = 0.832142857142857

# Row 6: Chemours Louisville Works
$ws.Range("A6").Value = "Chemours Louisville Works"
$ws.Range("B6").Value = "Louisville"
$ws.Range("C6").Value = 3707770
$ws.Range("D6").Value = 461
$ws.Range("E6").Value = 480.615153800134
$ws.Range("F6").Value = 553668
$ws.Range("G6").Value = 1151.99863263206
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0.069
$ws.Range("J6").Value = 387792
$ws.Range("K6").Value = 134546
$ws.Range("L6").Value = 1173
$ws.Range("M6").Value = 10148
$ws.Range("N6").Value = 30629
$ws.Range("O6").Value = 50.7676734234234
$ws.Range("P6").Value = 8.55059252947863
$ws.Range("Q6").Value = 10.5272202506423
$ws.Range("R6").Value = 30.1304347826087
$ws.Range("S6").Value = 0.414782608695652

# Row 7: Iofina Chemical Inc.
$ws.Range("A7").Value = "Iofina Chemical Inc."
$ws.Range("B7").Value = "Covington"
$ws.Range("C7").ClearContents()
$ws.Range("D7").Value = 414
$ws.Range("E7").Value = 482.413781082325
$ws.Range("F7").Value = 563753
$ws.Range("G7").Value = 1168.60882111449
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0.11
$ws.Range("J7").Value = 456893
$ws.Range("K7").Value = 74121
$ws.Range("L7").Value = 626
$ws.Range("M7").Value = 10942
$ws.Range("N7").Value = 18689
$ws.Range("O7").Value = 60.7478819095477
$ws.Range("P7").Value = 9.02237931054774
$ws.Range("Q7").Value = 9.70246182079913
$ws.Range("R7").Value = 29.636803874092
$ws.Range("S7").Value = 0.387409200968523

# Row 8: Honeywell International - Geismar Complex
$ws.Range("A8").Value = "Honeywell International - Geismar Complex"
$ws.Range("B8").Value = "Geismar"
$ws.Range("C8").Value = 413584
$ws.Range("D8").Value = 69
$ws.Range("E8").Value = 555.22662309248
$ws.Range("F8").Value = 182469
$ws.Range("G8").Value = 328.638779934022
$ws.Range("H8").Value = 1
$ws.Range("I8").Value = 0.52
$ws.Range("J8").Value = 119600
$ws.Range("K8").Value = 49943
$ws.Range("L8").Value = 141
$ws.Range("M8").Value = 3623
$ws.Range("N8").Value = 9220
$ws.Range("O8").Value = 72.2608507462687
$ws.Range("P8").Value = 5.92115845564032
$ws.Range("Q8").Value = 6.6311776732443
$ws.Range("R8").Value = 79.4202898550725
$ws.Range("S8").Value = 0.521739130434783

# Row 9: Mexichem Fluor Inc.
$ws.Range("A9").Value = "Mexichem Fluor Inc."
$ws.Range("B9").Value = "Saint Gabriel"
$ws.Range("C9").Value = 18331
$ws.Range("D9").Value = 68
$ws.Range("E9").Value = 506.724616424918
$ws.Range("F9").Value = 167063
$ws.Range("G9").Value = 329.691896909757
$ws.Range("H9").Value = 1
$ws.Range("I9").Value = 0.47
$ws.Range("J9").Value = 103097
$ws.Range("K9").Value = 51556
$ws.Range("L9").Value = 181
$ws.Range("M9").Value = 4141
$ws.Range("N9").Value = 8749
$ws.Range("O9").Value = 75.8387313432836
$ws.Range("P9").Value = 5.71423633158641
$ws.Range("Q9").Value = 7.25449505337132
$ws.Range("R9").Value = 80.5882352941177
$ws.Range("S9").Value = 0.557352941176471

# Row 10: Islechem LLC
$ws.Range("A10").Value = "Islechem LLC"
$ws.Range("B10").Value = "Grand Island"
$ws.Range("C10").ClearContents()
$ws.Range("D10").Value = 280
$ws.Range("E10").Value = 221.27861589355
$ws.Range("F10").Value = 303364
$ws.Range("G10").Value = 1370.95940687707
$ws.Range("H10").Value = 1
$ws.Range("I10").Value = 0.086
$ws.Range("J10").Value = 246489
$ws.Range("K10").Value = 27496
$ws.Range("L10").Value = 2358
$ws.Range("M10").Value = 9273
$ws.Range("N10").Value = 17628
$ws.Range("O10").Value = 57.1713522727273
$ws.Range("P10").Value = 8.43956350709922
$ws.Range("Q10").Value = 9.26747240535932
$ws.Range("R10").Value = 20.1428571428571
$ws.Range("S10").Value = 0.25

# Row 11: Chemours - Corpus Christi Plant
$ws.Range("A11").Value = "Chemours - Corpus Christi Plant"
$ws.Range("B11").Value = "Gregory"
$ws.Range("C11").Value = 17240
$ws.Range("D11").Value = 33
$ws.Range("E11").Value = 429.942996016745
$ws.Range("F11").Value = 51362
$ws.Range("G11").Value = 119.462348441186
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0.7
$ws.Range("J11").Value = 46644
$ws.Range("K11").Value = 1062
$ws.Range("L11").Value = 322
$ws.Range("M11").Value = 623
$ws.Range("N11").Value = 20610
$ws.Range("O11").Value = 62.32375
$ws.Range("P11").Value = 5.87796124304657
$ws.Range("Q11").Value = 5.93507974024599
$ws.Range("R11").Value = 20
$ws.Range("S11").Value = 0.212121212121212
